$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append new row 4 ---
$wsProximity = $wb.Worksheets.Item("Proximity")

$dateCellP = $wsProximity.Cells.Item(4, 1)
$dateCellP.NumberFormat = "@"
$dateCellP.Value = "2026-01-28"
$dateCellP.ClearFormats()

$wsProximity.Cells.Item(4, 2).Value = "18:46:33"
$wsProximity.Cells.Item(4, 3).Value = "18:00"
$wsProximity.Cells.Item(4, 4).Value = "Living Room Main Door"
$wsProximity.Cells.Item(4, 5).Value = "ENTER"
$wsProximity.Cells.Item(4, 6).Value = "User ENTERED Living Room Main Door"

# --- Camera sheet: append new row 4 ---
$wsCamera = $wb.Worksheets.Item("Camera")

$dateCellC = $wsCamera.Cells.Item(4, 1)
$dateCellC.NumberFormat = "@"
$dateCellC.Value = "2026-01-28"
$dateCellC.ClearFormats()

$wsCamera.Cells.Item(4, 2).Value = "18:46:34"
$wsCamera.Cells.Item(4, 3).Value = "18:00"
$wsCamera.Cells.Item(4, 4).Value = "Living Room Main Door"
$wsCamera.Cells.Item(4, 5).Value = "Image Captured"
$wsCamera.Cells.Item(4, 6).Value = "Active"

$wb.Save()
